$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to text
# (NumberFormat "@" during the write, then Style reset to Normal so the
# stored cell keeps the default/general style index like the rest of the sheet).

$ws.Range('D2').Value = '63.234.50'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '3.050.91'
$ws.Range('E3').Value = '  -2.93%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.65%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').Value = '3.046.50'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('E10').Value = '  -3.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.449'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.88'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('D16').Value = '3.556.79'
$ws.Range('E16').Value = '  -2.89%  '
$ws.Range('D17').Value = '63.309.54'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').Value = '3.050.53'
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '471.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.704'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('E29').Value = '  +2.01%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').Value = '0.0₃0821'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.20'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.70'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '439.21'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.287'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('E45').Value = '  +3.02%  '
$ws.Range('E46').Value = '  -4.08%  '
$ws.Range('D47').Value = '2.786.92'
$ws.Range('E47').Value = '  -3.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.77%  '
$ws.Range('E51').Value = '  -0.18%  '
